# Corrects several strings in the "en" language sheet (text tweaks / typo
# fixes) and updates the last active selection, matching the upstream
# "corrected texts, updated project to latest unity" commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("en")

$corrections = @{
    33  = "With these operators, we can get rid of them with ease."
    56  = "Before we proceed, let's learn some tricks with multiplication!"
    59  = "For example: 2 times 3 and 3 times 2 equal 6."
    60  = "With this trick, you only have to remember half of the multiplication table!"
    64  = "In multiples of three: double the number and then add the original number."
    65  = "For example, 3 times 6: double 6 to get 12, and then add 6 to get 18."
    67  = "For example, 4 times 6: double 6 to get 12, and then double 12 to get 24."
    69  = "Rearranging the equation and replacing division with multiplication can help."
    71  = "For the next mission, some blobs must be matched using division. Go ahead and try it out."
    78  = "For example, 6 times 7: multiply 7 by 5 to get 35, and then add 7 to get 42."
    88  = "For example, 8 times 4: double 4 to get 8, double 8 to get 16, and finally double 16 to get 32."
    90  = "When using multiplication, this is a good trick for splitting up a number, then computing each one separately."
    92  = "Now onward to the next mission!"
    95  = "For example, 9 times 6: multiply 6 by 10 to get 60, and then subtract 6 to get 54."
    98  = "Now it's time to clean up the last remaining blobs. Good luck!"
    100 = "You have banished all of the blobs! Earth is safe!"
}

foreach ($row in $corrections.Keys) {
    $ws.Cells.Item($row, 2).Value = $corrections[$row]
}

# Leave the same cell selected as in the saved workbook.
$ws.Range("B95").Select()
